$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.098.88'
$ws.Range('E2').Value = '  -0.28%  '
$ws.Range('D3').Value = '1.830.46'
$ws.Range('E3').Value = '  -0.01%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.009'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.45'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.23%  '
$ws.Range('E6').Value = '  -0.21%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4621'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.81%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3700'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.97%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07345'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.77%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8731'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.84%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07971'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.03%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '19.84'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.37%  '
$ws.Range('D13').Value = '1.784.26'
$ws.Range('E13').Value = '  -6.24%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.343'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.75%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.553'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.25%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '91.94'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.53%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.010'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.22%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008882'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.006'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.34%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.67'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.22%  '
$ws.Range('D21').Value = '26.913.73'
$ws.Range('E21').Value = '  -2.56%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.135'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.19%  '
$ws.Range('E23').Value = '  -0.30%  '
$ws.Range('D24').Value = '1.989.58'
$ws.Range('E24').Value = '  -4.69%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.46'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.843'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.10%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.64'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.62%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.081'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.34%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.096'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.72%  '
$ws.Range('E30').Value = '  -1.04%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08874'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.69%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.977'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.16%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7338'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.54%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.138'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.39%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.465'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.65%  '
$ws.Range('E37').Value = '  -1.49%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01946'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.24%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05242'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.22%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.945'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.15%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.149'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.59%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5189'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.45%  '
$ws.Range('E43').Value = '  -0.62%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8614'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -14.78%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.226'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.76%  '
$ws.Range('E46').Value = '  -1.17%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.27'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.28%  '
$ws.Range('B48').Value = 'PaxDollar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.008'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.13%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '102.32'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.09%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.628'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.63%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06230'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.88%  '
